# Applies the "Added many more features" revision to the Lucky Tree review.
#
# Strategy: every change in the target diff turns out to be a pure text
# substitution at a stable paragraph position - even the two bullet lists
# that look like insert/delete/reorder collapse to 1:1 text swaps once you
# line up old vs. new bullet order (the bullet counts are unchanged: 4 -> 4
# under "What we like", 2 -> 2 under "What we don't like"). So we just
# rewrite each paragraph's text in place by index, which keeps every other
# part of the document (styles, spacing, other paragraphs) untouched.

$d = $word.ActiveDocument

function Set-ParagraphText($paraIndex, $expectedOldText, $newText) {
    $p = $d.Paragraphs.Item($paraIndex)
    # Paragraph.Range.Text includes the trailing paragraph-mark character,
    # so trim it before comparing against the plain expected string.
    $old = $p.Range.Text.TrimEnd("`r", "`v", "`n")
    if ($old -ne $expectedOldText) {
        Write-Output ("WARNING: paragraph " + $paraIndex + " text was [" + $old + "], expected [" + $expectedOldText + "]")
    }
    $r = $d.Range($p.Range.Start, $p.Range.End)
    $r.Text = $newText
}

# --- 1 & 4: headline title, repeated verbatim further down as bold text ---
$d.Content.Find.Execute(
    "Play Lucky Tree for Free - A Unique Asian-Themed Slot",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Play Lucky Tree Free - Exciting Asian-themed Slot Game", 2) | Out-Null

# --- 2: "What we like" bullets ------------------------------------------------
# Before:                                              After:
#  36 Easy-to-understand gameplay and medium volatility  Easy to understand gameplay mechanics
#  37 Visually stunning graphics and symbols             Medium volatility for a good number of wins
#  38 Unique appeal in the Asian-themed slot category    Coins falling from the Lucky Tree for additional wins
#  39 Access to bonus mode with cat and dragon symbols   Visually stunning graphics and symbols
Set-ParagraphText 36 "Easy-to-understand gameplay and medium volatility" "Easy to understand gameplay mechanics"
Set-ParagraphText 37 "Visually stunning graphics and symbols" "Medium volatility for a good number of wins"
Set-ParagraphText 38 "Unique appeal in the Asian-themed slot category" "Coins falling from the Lucky Tree for additional wins"
Set-ParagraphText 39 "Access to bonus mode with cat and dragon symbols" "Visually stunning graphics and symbols"

# --- 3: "What we don't like" bullets -----------------------------------------
# Before:                                      After:
#  41 Sound department could be improved        Automatic spins are blocked during the coin falling feature
#  42 Lucky Tree feature blocks automatic spins  Sound department could be improved
Set-ParagraphText 41 "Sound department could be improved" "Automatic spins are blocked during the coin falling feature"
Set-ParagraphText 42 "Lucky Tree feature blocks automatic spins" "Sound department could be improved"

# --- 5: meta description (italic line) ---------------------------------------
$d.Content.Find.Execute(
    "Explore the nuanced theme of Asian culture while playing Lucky Tree, an online slot with visually stunning graphics and bonus modes. Play for free now.",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Read our review of Lucky Tree, a visually stunning Asian-themed slot game with free play.", 2) | Out-Null
